# Regenerate save_data "K" (strikeouts) column values for estrada_jeremiah.xlsx
# Per commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
# Only the values in column G (header "K") change for this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    9  = 4
    10 = 3
    11 = 1
    12 = 2
    13 = 3
    16 = 2
    17 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
